$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '59.636.04'
Set-TextValue $ws.Range('E2') '  -0.07%  '
Set-TextValue $ws.Range('D3') '2.398.50'
Set-TextValue $ws.Range('E3') '  -0.71%  '
Set-TextValue $ws.Range('E4') '  +0.19%  '
Set-TextValue $ws.Range('D5') '550.16'
Set-TextValue $ws.Range('E5') '  -0.64%  '
Set-TextValue $ws.Range('D6') '136.42'
Set-TextValue $ws.Range('E6') '  -0.78%  '
Set-TextValue $ws.Range('E7') '  +0.17%  '
Set-TextValue $ws.Range('D8') '0.590'
Set-TextValue $ws.Range('E8') '  +3.61%  '
Set-TextValue $ws.Range('E9') '  -2.04%  '
Set-TextValue $ws.Range('D10') '5.67'
Set-TextValue $ws.Range('E10') '  -1.96%  '
Set-TextValue $ws.Range('E11') '  -0.98%  '
Set-TextValue $ws.Range('D12') '0.353'
Set-TextValue $ws.Range('E12') '  -2.55%  '
Set-TextValue $ws.Range('D13') '25.34'
Set-TextValue $ws.Range('E13') '  +2.74%  '
Set-TextValue $ws.Range('D14') '2.822.85'
Set-TextValue $ws.Range('E14') '  -0.76%  '
Set-TextValue $ws.Range('D15') '59.629.01'
Set-TextValue $ws.Range('E15') '  +0.12%  '
Set-TextValue $ws.Range('E16') '  -2.08%  '
Set-TextValue $ws.Range('D17') '2.402.71'
Set-TextValue $ws.Range('E17') '  -0.80%  '
Set-TextValue $ws.Range('D18') '11.27'
Set-TextValue $ws.Range('E18') '  -0.53%  '
Set-TextValue $ws.Range('D19') '4.39'
Set-TextValue $ws.Range('E19') '  -1.56%  '
Set-TextValue $ws.Range('D20') '327.84'
Set-TextValue $ws.Range('E20') '  -2.19%  '
Set-TextValue $ws.Range('D21') '6.62'
Set-TextValue $ws.Range('E21') '  -4.65%  '
Set-TextValue $ws.Range('D22') '0.999'
Set-TextValue $ws.Range('E22') '  +0.00%  '
Set-TextValue $ws.Range('D23') '66.12'
Set-TextValue $ws.Range('E23') '  +2.46%  '
Set-TextValue $ws.Range('E24') '  +1.20%  '
Set-TextValue $ws.Range('D25') '8.60'
Set-TextValue $ws.Range('E25') '  -0.84%  '
Set-TextValue $ws.Range('D27') '1.36'
Set-TextValue $ws.Range('E27') '  -1.97%  '
Set-TextValue $ws.Range('D28') '0.0₃0767'
Set-TextValue $ws.Range('E28') '  -2.46%  '
Set-TextValue $ws.Range('D29') '1.76'
Set-TextValue $ws.Range('E29') '  -2.35%  '
Set-TextValue $ws.Range('D30') '168.49'
Set-TextValue $ws.Range('E30') '  -1.36%  '
Set-TextValue $ws.Range('D31') '6.05'
Set-TextValue $ws.Range('E31') '  -3.54%  '
Set-TextValue $ws.Range('D32') '18.55'
Set-TextValue $ws.Range('E32') '  -0.88%  '
Set-TextValue $ws.Range('E33') '  -2.23%  '
Set-TextValue $ws.Range('E34') '  -0.01%  '
Set-TextValue $ws.Range('E35') '  -1.56%  '
Set-TextValue $ws.Range('D36') '0.999'
Set-TextValue $ws.Range('E36') '  -0.18%  '
Set-TextValue $ws.Range('D37') '4.18'
Set-TextValue $ws.Range('E37') '  -2.54%  '
Set-TextValue $ws.Range('D38') '1.59'
Set-TextValue $ws.Range('E38') '  -1.96%  '
Set-TextValue $ws.Range('D39') '313.08'
Set-TextValue $ws.Range('E39') '  +2.25%  '
Set-TextValue $ws.Range('D40') '0.406'
Set-TextValue $ws.Range('E40') '  -2.95%  '
Set-TextValue $ws.Range('D41') '3.66'
Set-TextValue $ws.Range('E41') '  -2.57%  '
Set-TextValue $ws.Range('D42') '138.53'
Set-TextValue $ws.Range('E42') '  -2.73%  '
Set-TextValue $ws.Range('D43') '0.0967'
Set-TextValue $ws.Range('E43') '  +0.17%  '
Set-TextValue $ws.Range('D44') '0.0514'
Set-TextValue $ws.Range('E44') '  -1.81%  '
Set-TextValue $ws.Range('D45') '19.42'
Set-TextValue $ws.Range('E45') '  +1.50%  '
Set-TextValue $ws.Range('D46') '0.575'
Set-TextValue $ws.Range('E46') '  +0.58%  '
Set-TextValue $ws.Range('B47') 'VeChain'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D47') '0.0223'
Set-TextValue $ws.Range('E47') '  -1.57%  '
Set-TextValue $ws.Range('B48') 'Polygon'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D48') '0.387'
Set-TextValue $ws.Range('E48') '  -2.97%  '
Set-TextValue $ws.Range('D49') '17.57'
Set-TextValue $ws.Range('E49') '  -1.71%  '
Set-TextValue $ws.Range('D50') '11.07'
Set-TextValue $ws.Range('E50') '  +0.22%  '
Set-TextValue $ws.Range('D51') '1.56'
Set-TextValue $ws.Range('E51') '  -3.22%  '
